$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 3191
$ws.Range("B12").Value = 496

$ws.Range("B64").Value = 4591
$ws.Range("B65").Value = 4591
$ws.Range("B66").Value = 4367
$ws.Range("B67").Value = 4383
$ws.Range("B68").Value = 2929

$ws.Range("B80").Value = 5077
$ws.Range("B81").Value = 4243
$ws.Range("B82").Value = 3865

$ws.Range("B87").Value = 1534

$ws.Range("B90").Value = 206

$ws.Range("B95").Value = 385
$ws.Range("B96").Value = 385

$ws.Range("B107").Value = 128

$ws.Range("A110").Value = 17
$ws.Range("B110").Value = 60

$ws.Range("A111").Value = 162
$ws.Range("B111").Value = 60

$ws.Range("A112").Value = 158
$ws.Range("B112").Value = 3
